$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Tipo", matching the style of existing header cells (A1:C1)
$ws.Range("D1").Value = "Tipo"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats

# Update MSE (column B) and R2 (column C) values for rows 2-4
$newMSE = 0.1628417475493137
$newR2 = 0.9983795454146821

$ws.Range("B2").Value = $newMSE
$ws.Range("C2").Value = $newR2
$ws.Range("D2").Value = "single"

$ws.Range("B3").Value = $newMSE
$ws.Range("C3").Value = $newR2
$ws.Range("D3").Value = "single"

$ws.Range("B4").Value = $newMSE
$ws.Range("C4").Value = $newR2
$ws.Range("D4").Value = "single"
